$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.184.02"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "2.304.94"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.02"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.12"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.44%  "
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.516"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +5.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.32"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +9.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0793"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.68"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +11.66%  "
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("E14").Value = "  +3.68%  "
$ws.Range("D15").Value = "2.664.55"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").Value = "2.342.50"
$ws.Range("E16").Value = "  +3.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.804"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("D18").Value = "43.045.02"
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.70"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +11.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.21"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.00%  "
$ws.Range("E21").Value = "  +1.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.10"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.84"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.24"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +14.78%  "
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.45"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.06"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.20"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.89"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.88"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.18"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.04"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.84"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.50%  "
$ws.Range("E35").Value = "  -0.95%  "
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0696"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.91%  "
$ws.Range("E38").Value = "  +1.54%  "
$ws.Range("E39").Value = "  +3.93%  "
$ws.Range("E40").Value = "  +2.07%  "
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.32"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.98%  "
$ws.Range("D43").Value = "1.992.81"
$ws.Range("E43").Value = "  +2.02%  "
$ws.Range("E44").Value = "  +4.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.20"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +5.33%  "
$ws.Range("E46").Value = "  +1.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.90"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.78"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +6.53%  "
$ws.Range("E49").Value = "  +5.91%  "
$ws.Range("D50").Value = "2.530.87"
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("E51").Value = "  +1.90%  "
